# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" sheets to the newly generated values.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Map of row -> new F value for the "展览" sheet
$updates1 = @{
    3  = 3352
    5  = 2421
    9  = 1083
    14 = 96
    16 = 8456
    17 = 370
    19 = 249
    25 = 1150
    27 = 1986
    30 = 1730
    32 = 1914
    36 = 78
    37 = 179
    40 = 56
    42 = 402
    43 = 112
}

foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# Map of row -> new F value for the "全部类型" sheet
$updates4 = @{
    3  = 3352
    5  = 2421
    10 = 1083
    14 = 96
    16 = 8456
    17 = 370
    20 = 249
    26 = 1150
    28 = 1986
    30 = 1730
    32 = 1914
    36 = 78
    37 = 179
    40 = 56
    42 = 402
    47 = 112
}

foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
